# Refresh the cryptos price/volume snapshot (Price column D, Volume(1h) column E).
# Note: some Price values (e.g. "23.44") look numeric, so a leading apostrophe is
# used to force them to stay plain text, matching how the sheet already stores
# every Price/Volume cell (the column already mixes plain decimals with
# "thousands.dotted" values like "1.649.06", so everything here must remain text).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.997.90'
$ws.Range('E2').Value = '  +1.81%  '
$ws.Range('D3').Value = '1.649.06'
$ws.Range('E3').Value = '  +1.94%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').Value = '''213.76'
$ws.Range('E5').Value = '  +1.44%  '
$ws.Range('E6').Value = '  +0.55%  '
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('D8').Value = '''23.44'
$ws.Range('E8').Value = '  +2.82%  '
$ws.Range('D9').Value = '''0.266'
$ws.Range('E9').Value = '  +1.67%  '
$ws.Range('E10').Value = '  +0.27%  '
$ws.Range('E11').Value = '  -1.50%  '
$ws.Range('D12').Value = '1.880.30'
$ws.Range('E12').Value = '  +1.78%  '
$ws.Range('D13').Value = '1.652.22'
$ws.Range('E13').Value = '  +2.21%  '
$ws.Range('E14').Value = '  +1.17%  '
$ws.Range('D15').Value = '''0.566'
$ws.Range('E15').Value = '  +3.03%  '
$ws.Range('D16').Value = '''65.64'
$ws.Range('E16').Value = '  +0.84%  '
$ws.Range('D17').Value = '27.993.60'
$ws.Range('E17').Value = '  +1.88%  '
$ws.Range('D18').Value = '''233.12'
$ws.Range('E18').Value = '  +0.98%  '
$ws.Range('D19').Value = '''7.70'
$ws.Range('E19').Value = '  +2.59%  '
$ws.Range('E20').Value = '  +0.62%  '
$ws.Range('E21').Value = '  -0.02%  '
$ws.Range('D22').Value = '''10.69'
$ws.Range('E22').Value = '  +5.18%  '
$ws.Range('E23').Value = '  +2.81%  '
$ws.Range('D24').Value = '''2.15'
$ws.Range('E24').Value = '  +4.01%  '
$ws.Range('D25').Value = '''152.23'
$ws.Range('E25').Value = '  +0.95%  '
$ws.Range('D26').Value = '''6.93'
$ws.Range('E26').Value = '  +1.44%  '
$ws.Range('D27').Value = '''15.79'
$ws.Range('E27').Value = '  +1.71%  '
$ws.Range('D28').Value = '''0.112'
$ws.Range('E28').Value = '  +0.19%  '
$ws.Range('E29').Value = '  -0.03%  '
$ws.Range('E30').Value = '  +1.55%  '
$ws.Range('E31').Value = '  +0.35%  '
$ws.Range('E32').Value = '  +2.79%  '
$ws.Range('D33').Value = '1.445.88'
$ws.Range('E33').Value = '  -1.56%  '
$ws.Range('D34').Value = '''3.09'
$ws.Range('E34').Value = '  +0.42%  '
$ws.Range('E35').Value = '  +2.23%  '
$ws.Range('D36').Value = '''2.33'
$ws.Range('E36').Value = '  -0.45%  '
$ws.Range('D37').Value = '''0.890'
$ws.Range('E37').Value = '  +3.55%  '
$ws.Range('E38').Value = '  +1.18%  '
$ws.Range('E39').Value = '  +0.42%  '
$ws.Range('D40').Value = '''0.921'
$ws.Range('E40').Value = '  -3.15%  '
$ws.Range('D41').Value = '''69.37'
$ws.Range('E41').Value = '  +2.13%  '
$ws.Range('E42').Value = '  +3.56%  '
$ws.Range('E44').Value = '  -0.26%  '
$ws.Range('E45').Value = '  +1.15%  '
$ws.Range('D46').Value = '''5.41'
$ws.Range('E46').Value = '  +2.85%  '
$ws.Range('D47').Value = '''1.80'
$ws.Range('E47').Value = '  +5.27%  '
$ws.Range('D48').Value = '1.789.40'
$ws.Range('E48').Value = '  +1.62%  '
$ws.Range('D49').Value = '''89.00'
$ws.Range('E49').Value = '  +2.87%  '
$ws.Range('E50').Value = '  -0.16%  '
$ws.Range('E51').Value = '  +0.29%  '
